$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were deleted entirely in the fixed export (naive forecaster bug)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update recalculated forecast values (tiny floating point deltas from the bugfix)
$ws.Range("E3").Value = 9.591339540850829
$ws.Range("E4").Value = 4.422525088127305
$ws.Range("C6").Value = -14.45332333832744
$ws.Range("E6").Value = -2.928447329610051
$ws.Range("E7").Value = -2.225127715916664
$ws.Range("C8").Value = 8.600536527919612
$ws.Range("C9").Value = 9.399485634179205
$ws.Range("C11").Value = 5.169490031659651
$ws.Range("E11").Value = 9.213376886330327
$ws.Range("C12").Value = 4.639893381363192
$ws.Range("C13").Value = -0.3722371047999995
$ws.Range("E13").Value = 2.684220738731979
$ws.Range("E14").Value = 2.429116709932599
$ws.Range("C15").Value = 4.098801479368319
$ws.Range("E16").Value = 3.941300050092877
$ws.Range("E17").Value = 2.714258593289998
$ws.Range("C19").Value = 2.352205130086094
$ws.Range("C21").Value = 4.083548352538391
$ws.Range("E21").Value = 3.58625614607444
$ws.Range("C22").Value = 4.695933104194361
$ws.Range("C24").Value = 4.861590900330715
$ws.Range("C25").Value = 5.402237127943765
$ws.Range("E25").Value = 4.104053120889195
$ws.Range("C29").Value = 0.8513583007189407
$ws.Range("E29").Value = 2.225279621195853
$ws.Range("C31").Value = 1.015697339178057
$ws.Range("E31").Value = 2.122104735451624
$ws.Range("E32").Value = -0.6322362079330235
$ws.Range("E34").Value = -1.352810423674367
$ws.Range("C35").Value = 4.074459326939817
$ws.Range("E35").Value = -0.2414327668618488
$ws.Range("E36").Value = 1.077755602068309
$ws.Range("E37").Value = 1.148476797857967
$ws.Range("E38").Value = -0.3934198590721305
$ws.Range("E41").Value = 1.534407168230811
$ws.Range("C42").Value = 5.120680133083622
$ws.Range("C43").Value = 5.356482122456163
$ws.Range("E43").Value = 12.6296844023545
$ws.Range("C46").Value = -0.5532735011319123
$ws.Range("C47").Value = -2.464475897442031
$ws.Range("E48").Value = 1.793234865396331
$ws.Range("C49").Value = -0.8995735674421024
$ws.Range("E49").Value = 0.3402056885013494
$ws.Range("E50").Value = -1.28528149926006
$ws.Range("E51").Value = -1.655020334777801
$ws.Range("E52").Value = -0.6714033493142035
$ws.Range("E53").Value = -0.5945514555738662

Write-Output "applied naive forecaster bugfix values"
